$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'23.922.77"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.22%  "

$ws.Range("D3").Value = "'1.655.71"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.01%  "

$ws.Range("D4").Value = "'1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'308.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.78%  "

$ws.Range("E6").Value = "  -0.36%  "

$ws.Range("D7").Value = "'0.3884"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.51%  "

$ws.Range("D8").Value = "'0.3844"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.50%  "

$ws.Range("D9").Value = "'51.63"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.04%  "

$ws.Range("D10").Value = "'1.358"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.24%  "

$ws.Range("E11").Value = "  -0.13%  "

$ws.Range("D12").Value = "'0.08469"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.67%  "

$ws.Range("D13").Value = "'23.84"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.22%  "

$ws.Range("D14").Value = "'7.191"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.12%  "

$ws.Range("D15").Value = "'7.967"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.40%  "

$ws.Range("D16").Value = "'0.00001306"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.61%  "

$ws.Range("D17").Value = "'1.657.39"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.39%  "

$ws.Range("D18").Value = "'94.21"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.45%  "

$ws.Range("D19").Value = "'0.06982"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.83%  "

$ws.Range("D20").Value = "'19.81"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.19%  "

$ws.Range("D21").Value = "'6.954"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.20%  "

$ws.Range("D22").Value = "'1.002"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.29%  "

$ws.Range("D23").Value = "'13.60"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.78%  "

$ws.Range("D24").Value = "'23.928.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.23%  "

$ws.Range("D25").Value = "'2.495"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.24%  "

$ws.Range("D26").Value = "'3.085"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +9.45%  "

$ws.Range("D27").Value = "'22.17"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.59%  "

$ws.Range("D28").Value = "'152.96"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.12%  "

$ws.Range("D29").Value = "'139.81"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.51%  "

$ws.Range("D30").Value = "'5.308"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.14%  "

$ws.Range("D31").Value = "'7.867"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.57%  "

$ws.Range("D32").Value = "'2.482"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.02%  "

$ws.Range("D33").Value = "'1.841.15"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.29%  "

$ws.Range("D34").Value = "'1.034"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.57%  "

$ws.Range("D35").Value = "'0.08097"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.85%  "

$ws.Range("D36").Value = "'0.02974"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.44%  "

$ws.Range("D37").Value = "'10.99"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.34%  "

$ws.Range("D38").Value = "'6.691"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.76%  "

$ws.Range("D39").Value = "'0.2691"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.57%  "

$ws.Range("D40").Value = "'0.09135"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.20%  "

$ws.Range("D41").Value = "'13.55"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.11%  "

$ws.Range("D42").Value = "'0.7519"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.01%  "

$ws.Range("D43").Value = "'1.417"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.18%  "

$ws.Range("D44").Value = "'16.47"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.52%  "

$ws.Range("D45").Value = "'0.6965"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.02%  "

$ws.Range("D46").Value = "'2.472"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.47%  "

$ws.Range("D47").Value = "'4.078"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.64%  "

$ws.Range("D48").Value = "'1.001"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.31%  "

$ws.Range("D49").Value = "'0.08281"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.48%  "

$ws.Range("D50").Value = "'135.41"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.10%  "

$ws.Range("D51").Value = "'1.230"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.52%  "
